$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.490.61"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "3.015.26"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.523"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").Value = "3.011.49"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "3.518.92"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.501.14"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "3.016.75"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.692"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -6.44%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "0.0₃0791"
$ws.Range("E36").Value = "  -5.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "421.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.112"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "2.799.91"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("E51").Value = "  -1.35%  "
